$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = "30.038.45"
$ws.Range("E2").Value2 = "  -0.30%  "
$ws.Range("D3").Value2 = "1.869.41"
$ws.Range("E3").Value2 = "  -3.00%  "
$ws.Range("E4").Value2 = "  +0.24%  "
$ws.Range("D5").Value2 = "'319.06"
$ws.Range("E5").Value2 = "  -0.85%  "
$ws.Range("D6").Value2 = "'1.002"
$ws.Range("E6").Value2 = "  +0.20%  "
$ws.Range("D7").Value2 = "'0.5073"
$ws.Range("E7").Value2 = "  -1.84%  "
$ws.Range("D8").Value2 = "'0.3936"
$ws.Range("E8").Value2 = "  -1.52%  "
$ws.Range("D9").Value2 = "'0.08193"
$ws.Range("E9").Value2 = "  -3.49%  "
$ws.Range("D10").Value2 = "'42.11"
$ws.Range("E10").Value2 = "  -2.09%  "
$ws.Range("D11").Value2 = "'1.091"
$ws.Range("E11").Value2 = "  -3.01%  "
$ws.Range("D12").Value2 = "'22.65"
$ws.Range("E12").Value2 = "  +6.59%  "
$ws.Range("D13").Value2 = "1.868.60"
$ws.Range("E13").Value2 = "  -3.00%  "
$ws.Range("D14").Value2 = "'6.256"
$ws.Range("E14").Value2 = "  -1.24%  "
$ws.Range("D15").Value2 = "'7.148"
$ws.Range("E15").Value2 = "  -3.22%  "
$ws.Range("D16").Value2 = "'1.004"
$ws.Range("E16").Value2 = "  +0.27%  "
$ws.Range("D17").Value2 = "'92.13"
$ws.Range("E17").Value2 = "  -2.48%  "
$ws.Range("D18").Value2 = "'0.00001079"
$ws.Range("E18").Value2 = "  -3.47%  "
$ws.Range("D19").Value2 = "'0.06330"
$ws.Range("E19").Value2 = "  -6.39%  "
$ws.Range("D20").Value2 = "'17.83"
$ws.Range("E20").Value2 = "  -0.96%  "
$ws.Range("D21").Value2 = "'1.002"
$ws.Range("E21").Value2 = "  +0.21%  "
$ws.Range("D22").Value2 = "30.018.46"
$ws.Range("E22").Value2 = "  -0.40%  "
$ws.Range("D23").Value2 = "'5.798"
$ws.Range("E23").Value2 = "  -4.65%  "
$ws.Range("D24").Value2 = "'11.04"
$ws.Range("E24").Value2 = "  -1.71%  "
$ws.Range("E25").Value2 = "  +0.07%  "
$ws.Range("D26").Value2 = "2.088.40"
$ws.Range("E26").Value2 = "  -2.63%  "
$ws.Range("D27").Value2 = "'161.38"
$ws.Range("E27").Value2 = "  +0.97%  "
$ws.Range("D28").Value2 = "'20.93"
$ws.Range("E28").Value2 = "  -0.79%  "
$ws.Range("D29").Value2 = "'2.245"
$ws.Range("E29").Value2 = "  -9.25%  "
$ws.Range("D30").Value2 = "'126.70"
$ws.Range("E30").Value2 = "  -2.27%  "
$ws.Range("B31").Value2 = "Stellar"
$ws.Range("C31").Value2 = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D31").Value2 = "'0.1033"
$ws.Range("E31").Value2 = "  -2.44%  "
$ws.Range("B32").Value2 = "ImmutableX"
$ws.Range("C32").Value2 = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D32").Value2 = "'1.043"
$ws.Range("E32").Value2 = "  -3.33%  "
$ws.Range("D33").Value2 = "'5.864"
$ws.Range("E33").Value2 = "  -4.07%  "
$ws.Range("D34").Value2 = "'3.724"
$ws.Range("E34").Value2 = "  +1.52%  "
$ws.Range("D35").Value2 = "'0.02429"
$ws.Range("E35").Value2 = "  -2.93%  "
$ws.Range("D36").Value2 = "'5.208"
$ws.Range("E36").Value2 = "  +0.10%  "
$ws.Range("D37").Value2 = "'0.06344"
$ws.Range("E37").Value2 = "  -4.24%  "
$ws.Range("D38").Value2 = "'0.2139"
$ws.Range("E38").Value2 = "  -3.36%  "
$ws.Range("D39").Value2 = "'1.170"
$ws.Range("E39").Value2 = "  -6.23%  "
$ws.Range("D40").Value2 = "'8.515"
$ws.Range("E40").Value2 = "  -5.87%  "
$ws.Range("D41").Value2 = "'0.6275"
$ws.Range("E41").Value2 = "  -3.96%  "
$ws.Range("D42").Value2 = "'1.211"
$ws.Range("E42").Value2 = "  -2.34%  "
$ws.Range("D43").Value2 = "'11.26"
$ws.Range("E43").Value2 = "  -1.39%  "
$ws.Range("B44").Value2 = "Frax"
$ws.Range("C44").Value2 = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D44").Value2 = "'1.001"
$ws.Range("E44").Value2 = "  +0.13%  "
$ws.Range("B45").Value2 = "Decentraland"
$ws.Range("C45").Value2 = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D45").Value2 = "'0.5892"
$ws.Range("E45").Value2 = "  -4.25%  "
$ws.Range("B46").Value2 = "EnergySwap"
$ws.Range("C46").Value2 = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").Value2 = "'12.88"
$ws.Range("E46").Value2 = "  -2.07%  "
$ws.Range("B47").Value2 = "PancakeSwap"
$ws.Range("C47").Value2 = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D47").Value2 = "'3.638"
$ws.Range("E47").Value2 = "  -2.32%  "
$ws.Range("B48").Value2 = "NEARProtocol"
$ws.Range("C48").Value2 = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D48").Value2 = "'1.989"
$ws.Range("E48").Value2 = "  -3.43%  "
$ws.Range("B49").Value2 = "EOS"
$ws.Range("C49").Value2 = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D49").Value2 = "'1.206"
$ws.Range("E49").Value2 = "  -2.96%  "
$ws.Range("B50").Value2 = "Quant"
$ws.Range("C50").Value2 = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D50").Value2 = "'121.75"
$ws.Range("E50").Value2 = "  -3.03%  "
$ws.Range("B51").Value2 = "WEMIXTOKEN"
$ws.Range("C51").Value2 = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D51").Value2 = "'1.119"
$ws.Range("E51").Value2 = "  -2.75%  "
